# Auto-generated: updates cryptos Price (D) and Volume(1h) (E) columns
# per the commit diff for Sun Oct  1 04:46:36 UTC 2023 data refresh.
# Numeric-looking D-column price strings (e.g. "215.40") are forced to
# remain plain text (matching the source inlineStr cells) by switching the
# cell to a text NumberFormat before the write, then resetting the style
# back to Normal afterwards so no stray style/format residue is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.071.08'
$ws.Cells.Item(2, 5).Value = '  +0.46%  '

$ws.Cells.Item(3, 4).Value = '1.676.55'
$ws.Cells.Item(3, 5).Value = '  +0.31%  '

$ws.Cells.Item(4, 5).Value = '  +0.13%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '215.40'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.28%  '

$ws.Cells.Item(7, 5).Value = '  +0.08%  '

$ws.Cells.Item(8, 5).Value = '  +1.82%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '21.31'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +5.59%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0885'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.64%  '

$ws.Cells.Item(12, 4).Value = '1.912.24'

$ws.Cells.Item(13, 4).Value = '1.686.53'
$ws.Cells.Item(13, 5).Value = '  -0.23%  '

$ws.Cells.Item(14, 5).Value = '  +0.88%  '

$ws.Cells.Item(15, 5).Value = '  +1.50%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '65.99'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.72%  '

$ws.Cells.Item(17, 4).Value = '27.057.27'
$ws.Cells.Item(17, 5).Value = '  +0.41%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '237.51'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +1.46%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '8.13'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +1.10%  '

$ws.Cells.Item(20, 5).Value = '  +1.16%  '

$ws.Cells.Item(21, 5).Value = '  +0.06%  '

$ws.Cells.Item(22, 5).Value = '  +0.76%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '9.33'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +1.73%  '

$ws.Cells.Item(24, 5).Value = '  -2.23%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '146.79'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.61%  '

$ws.Cells.Item(26, 5).Value = '  +1.28%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '16.34'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +2.55%  '

$ws.Cells.Item(28, 5).Value = '  +0.31%  '

$ws.Cells.Item(29, 5).Value = '  +0.23%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.0498'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.03%  '

$ws.Cells.Item(31, 5).Value = '  +0.05%  '

$ws.Cells.Item(33, 4).Value = '1.545.65'
$ws.Cells.Item(33, 5).Value = '  +5.85%  '

$ws.Cells.Item(34, 5).Value = '  +1.61%  '

$ws.Cells.Item(35, 5).Value = '  +2.56%  '

$ws.Cells.Item(36, 5).Value = '  +3.26%  '

$ws.Cells.Item(37, 5).Value = '  -1.11%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.924'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +2.83%  '

$ws.Cells.Item(39, 5).Value = '  +2.01%  '

$ws.Cells.Item(40, 5).Value = '  +1.87%  '

$ws.Cells.Item(41, 5).Value = '  +0.08%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '67.64'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +1.72%  '

$ws.Cells.Item(43, 5).Value = '  -3.36%  '

$ws.Cells.Item(44, 5).Value = '  -2.05%  '

$ws.Cells.Item(45, 4).Value = '1.821.51'
$ws.Cells.Item(45, 5).Value = '  +0.67%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.781'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.06%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '90.68'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.07%  '

$ws.Cells.Item(48, 5).Value = '  +1.79%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.57'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +2.21%  '

$ws.Cells.Item(50, 5).Value = '  +2.53%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '8.03'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +4.89%  '
